# "Fortran Port Tasklist" update:
#  - bump the "Current as of:" date
#  - flip task #4 from Upcoming -> Paused and swap in its real description/note
#  - append tasks #5 and #6 (with #6 split into sub-steps 6, 6.01, 6.02, 6.03)
#
# NOTE on write order: this engine builds xl/sharedStrings.xml in first-seen
# order as cells are written, and the target file's table has
# "Allocate and initialize variables" registered before "Initialize simulation"
# even though "Initialize simulation" (row 9, C9) sits above "Allocate and
# initialize variables" (row 10, C10) in the sheet. To reproduce that table
# order exactly, C10 is written before C9 below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Current as of: 2021-10-21 -> 2021-11-17 (serial 44490 -> 44517)
$ws.Range("B1").Value = 44517

# Row 7 / Task #4: Upcoming -> Paused, new description + note, completion 0 -> 75%
$ws.Range("B7").Value = "Paused"
$ws.Range("C7").Value = "Create variables, structures, classes, etc."
$ws.Range("D7").Value = 0.75
$ws.Range("E7").Value = "Most important variables created in simple arrays instead of complex structures. Variables will be defined as necessary."

# Row 8 / Task #5: In Progress - Create configuration file
$ws.Range("A8").Value = 5
$ws.Range("A8").HorizontalAlignment = $ws.Range("A7").HorizontalAlignment
$ws.Range("B8").Value = "In Progress"
$ws.Range("C8").Value = "Create configuration file"
$ws.Range("D8").Value = 0
$ws.Range("D8").NumberFormat = $ws.Range("D7").NumberFormat
$ws.Range("E8").Value = "Localize simulation configuration to one place, instead of having to check multiple places before running (like in the C++ code)"

# Row 9 / Task #6: Future - Initialize simulation
$ws.Range("A9").Value = 6
$ws.Range("A9").HorizontalAlignment = $ws.Range("A7").HorizontalAlignment
$ws.Range("B9").Value = "Future"

# Row 10 / Task #6.01: Future - Allocate and initialize variables (written before C9, see NOTE above)
$ws.Range("A10").Value = 6.01
$ws.Range("A10").HorizontalAlignment = $ws.Range("A7").HorizontalAlignment
$ws.Range("B10").Value = "Future"
$ws.Range("C10").Value = "Allocate and initialize variables"

$ws.Range("C9").Value = "Initialize simulation"

# Row 11 / Task #6.02: Future - Enforce boundary conditions on initial conditions
$ws.Range("A11").Value = 6.02
$ws.Range("A11").HorizontalAlignment = $ws.Range("A7").HorizontalAlignment
$ws.Range("B11").Value = "Future"
$ws.Range("C11").Value = "Enforce boundary conditions on initial conditions"

# Row 12 / Task #6.03: Future - Calculate flow derivatives
$ws.Range("A12").Value = 6.03
$ws.Range("A12").HorizontalAlignment = $ws.Range("A7").HorizontalAlignment
$ws.Range("B12").Value = "Future"
$ws.Range("C12").Value = "Calculate flow derivatives"

# Selection moves to B2
$ws.Range("B2").Select()
